$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.720.66"
$ws.Range("E2").Value = '  -3.25%  '

$ws.Range("D3").Value = "'1.736.05"
$ws.Range("E3").Value = '  -5.65%  '

$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'238.72"
$ws.Range("E5").Value = '  -7.76%  '

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").Value = "'0.4899"
$ws.Range("E7").Value = '  -7.01%  '

$ws.Range("E8").Value = '  -7.09%  '

$ws.Range("D9").Value = "'0.2425"
$ws.Range("E9").Value = '  -23.02%  '

$ws.Range("D10").Value = "'0.06010"
$ws.Range("E10").Value = '  -11.54%  '

$ws.Range("D11").Value = "'1.740.40"
$ws.Range("E11").Value = '  -5.56%  '

$ws.Range("D12").Value = "'0.06695"
$ws.Range("E12").Value = '  -13.65%  '

$ws.Range("D13").Value = "'14.87"
$ws.Range("E13").Value = '  -20.34%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = "'0.5908"
$ws.Range("E14").Value = '  -24.11%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = "'76.80"
$ws.Range("E15").Value = '  -12.51%  '

$ws.Range("D16").Value = "'4.387"
$ws.Range("E16").Value = '  -12.31%  '

$ws.Range("D17").Value = "'0.9990"
$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("D19").Value = "'25.755.06"
$ws.Range("E19").Value = '  -3.15%  '

$ws.Range("D20").Value = "'11.39"
$ws.Range("E20").Value = '  -17.55%  '

$ws.Range("D21").Value = "'0.000006321"
$ws.Range("E21").Value = '  -20.13%  '

$ws.Range("D22").Value = "'1.960.65"
$ws.Range("E22").Value = '  -5.84%  '

$ws.Range("D23").Value = "'3.909"
$ws.Range("E23").Value = '  -14.90%  '

$ws.Range("D24").Value = "'5.118"
$ws.Range("E24").Value = '  -14.16%  '

$ws.Range("D25").Value = "'7.895"
$ws.Range("E25").Value = '  -15.22%  '

$ws.Range("D26").Value = "'135.60"
$ws.Range("E26").Value = '  -5.26%  '

$ws.Range("D27").Value = "'1.848"
$ws.Range("E27").Value = '  -16.47%  '

$ws.Range("D28").Value = "'1.427"
$ws.Range("E28").Value = '  -15.34%  '

$ws.Range("D29").Value = "'14.31"
$ws.Range("E29").Value = '  -15.59%  '

$ws.Range("D30").Value = "'100.47"
$ws.Range("E30").Value = '  -9.44%  '

$ws.Range("D31").Value = "'0.08156"
$ws.Range("E31").Value = '  -6.44%  '

$ws.Range("D32").Value = "'3.640"
$ws.Range("E32").Value = '  -12.84%  '

$ws.Range("D33").Value = "'3.289"
$ws.Range("E33").Value = '  -18.97%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = "'0.04357"
$ws.Range("E34").Value = '  -10.47%  '

$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").Value = "'2.682"
$ws.Range("E36").Value = '  -6.19%  '

$ws.Range("D37").Value = "'1.026"
$ws.Range("E37").Value = '  -9.79%  '

$ws.Range("D38").Value = "'0.6084"
$ws.Range("E38").Value = '  -16.75%  '

$ws.Range("D39").Value = "'2.780"
$ws.Range("E39").Value = '  -9.90%  '

$ws.Range("D40").Value = "'2.076"
$ws.Range("E40").Value = '  -7.38%  '

$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = "'101.98"
$ws.Range("E42").Value = '  -7.09%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = "'0.01489"
$ws.Range("E43").Value = '  -13.27%  '

$ws.Range("D44").Value = "'0.7919"
$ws.Range("E44").Value = '  -11.44%  '

$ws.Range("D45").Value = "'0.3802"
$ws.Range("E45").Value = '  -20.82%  '

$ws.Range("D46").Value = "'5.128"
$ws.Range("E46").Value = '  -13.36%  '

$ws.Range("D47").Value = "'6.077"
$ws.Range("E47").Value = '  -20.37%  '

$ws.Range("D48").Value = "'0.05081"
$ws.Range("E48").Value = '  -12.62%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'1.249"
$ws.Range("E49").Value = '  -11.47%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = "'29.66"
$ws.Range("E50").Value = '  -14.53%  '

$ws.Range("D51").Value = "'51.94"
$ws.Range("E51").Value = '  -13.12%  '
